$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet2: replace the lone sample value with the real "counter" report
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("A1").Value = "Schedule Number Counter"
$ws2.Range("B1").Value = 129
$ws2.Range("A2").Value = "Data Recorder Index"
$ws2.Range("B2").Value = 12

$ws2.Columns.Item(1).AutoFit() | Out-Null

# ---------------------------------------------------------------------
# Sheet3: brand-new worksheet with the CO/MO/Schedule/Delivery table
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "Sheet3"

$ws3.Range("A1").Value = "CONumber"
$ws3.Range("B1").Value = "MONumber"
$ws3.Range("C1").Value = "ScheduleNumber"
$ws3.Range("D1").Value = "DeliveryNumber"

# CONumber / MONumber / DeliveryNumber are long numeric-looking codes that
# must land in the sheet as TEXT (shared strings), not numbers - otherwise
# the leading/scale semantics of these reference numbers would change.
# ScheduleNumber (column C) is a genuine number.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$rows = @(
    @("3013691711", "1000004200", 21708102017, "13141381"),
    @("3013691714", "1000004203", 22008102017, "13141387"),
    @("3013691715", "1000004204", 22108102017, "13141389"),
    @("3013691716", "1000004205", 22208102017, "13141391"),
    @("3013691717", "1000004206", 22308102017, "13141393"),
    @("3013691718", "1000004207", 22408102017, "13141395"),
    @("3013691730", "1000004218", 22608102017, "13141496"),
    @("3013691731", "1000004219", 22708102017, "13141498"),
    @("3013691732", "1000004220", 22808102017, "13141500"),
    @("3013691733", "1000004221", 22908102017, "13141502")
)

$r = 2
foreach ($row in $rows) {
    Set-TextValue $ws3.Cells.Item($r, 1) $row[0]
    Set-TextValue $ws3.Cells.Item($r, 2) $row[1]
    $ws3.Cells.Item($r, 3).Value = $row[2]
    Set-TextValue $ws3.Cells.Item($r, 4) $row[3]
    $r = $r + 1
}

$ws3.Columns.Item(1).AutoFit() | Out-Null
$ws3.Columns.Item(2).AutoFit() | Out-Null
$ws3.Columns.Item(3).AutoFit() | Out-Null
$ws3.Columns.Item(4).AutoFit() | Out-Null

$ws3.Range("A2:D3").Select() | Out-Null

# Leave Sheet2 as the active/selected tab (matches the original workbook,
# which only gained a new trailing sheet).
$ws2.Activate()
$ws2.Range("B2").Select() | Out-Null
